$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

$ws.Cells.Item(14, 1).Value = "2025-08-28T18:31:57.556918"
$ws.Cells.Item(14, 2).Value = 13
$ws.Cells.Item(14, 3).Value = "全案件リスト"
$ws.Cells.Item(14, 4).Value = 53.8
$ws.Cells.Item(14, 5).Value = 6
$ws.Cells.Item(14, 6).Value = 4
$ws.Cells.Item(14, 7).Value = 13
